# Insert a new data row at row 168 (pushes existing rows 168:261 down to 169:262)
# and populate it with the new daily price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(168).Insert()

$ws.Range("A168").Value = 10
$ws.Range("B168").Value = "Vega Modelo de Temuco"
$ws.Range("C168").Value = "La Araucanía"
$ws.Range("D168").Value = 44777
$ws.Range("E168").Value = 9
$ws.Range("F168").Value = 100112043
$ws.Range("G168").Value = "Pepino dulce"
$ws.Range("H168").Value = "Cultivar IV Región"
$ws.Range("I168").Value = "Primera"
$ws.Range("J168").Value = 65
$ws.Range("K168").Value = 20000
$ws.Range("L168").Value = 20000
$ws.Range("M168").Value = 20000
$ws.Range("N168").Value = "$/bandeja 18 kilos"
$ws.Range("O168").Value = "Provincia de Limarí"
$ws.Range("P168").Value = 1111
$ws.Range("Q168").Value = 18
$ws.Range("R168").Value = "Hortaliza"
